# "Just adding finishing touches" - append five more skip-trace records to
# Sheet4, tidy up the stray date-format left on B3, and resize the columns
# so the new, wider data is readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the leftover date style on B3 (should just be a plain number,
#     same as B2) --------------------------------------------------------
$ws.Range("B3").Style = "Normal"

# --- new rows ------------------------------------------------------------
$ws.Range("A4").Value = "Michelle Semexant"
$ws.Range("B4").Value = 32287
$ws.Range("C4").Value = "4233836550 6786434050 6786434050 6783349534 6786434050 7708743159 7709891035 7709891035 6783343495 6783349534 6783349534 7707329164 6786434050 "
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = "At 2731 VALLEY GREEN DR from Thu Mar 21 2019 to Tue Aug 13 2019"

$ws.Range("A5").Value = "Henry Turner"
$ws.Range("B5").Value = 30923
$ws.Range("C5").Value = "2158372470 4047290980 6786343071 6789061865 8505904243 2085202431 6264754955 8502944299 8502944299 8502944299 8503850412 4042540441 4045561004 8503052955 8506560525 "
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "At 515 ROSEDOWN WAY since Thu Oct 05 2017"

$ws.Range("A6").Value = "Jason Hyman"
$ws.Range("B6").Value = "Fri Oct 18 1974"
$ws.Range("C6").Value = "4042293909 4046643327 4048738310 4048738310 6783586414 4042550808 4042550808 4043925135 4043925135 6783586414 6783586414 6785954385 6785954386 7709771941 5124979968 4046643327 "
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "At 590 BRIDGEWATER DR since Mon Jul 02 2012"

$ws.Range("A7").Value = "David Earley"
$ws.Range("B7").Value = "Sat Oct 28 1989"
$ws.Range("C7").Value = "4436945686 4436943766 5087985592 5087985592 7635444403 6123452837 6126707125 6129913083 7634584476 2147558136 2147558136 2147558136 "
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "At 5147 WELLSLEY BND from Sat Apr 02 2016 to Sun Nov 04 2018"

$ws.Range("A8").Value = "Curran Sullivan"
$ws.Range("B8").Value = "Wed Dec 19 1990"
$ws.Range("C8").Value = "8053251450 8053251036 8053251450 8057322769 8649334047 8053251450 8056887004 8056887004 8056889053 8057322769 8053251036 8053251037 "
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "At 2309 FORREST WALK since Mon Aug 12 2019"

# --- resize columns A-C so the new (longer) names / phone-number blob fit
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).ColumnWidth = 66.3

# --- dimension / selection follow the newly-used range --------------------
$ws.Range("C12").Select()

$wb.Save()
